$wb = $excel.ActiveWorkbook

# --- Sheet: Heat Generators ---
$ws1 = $wb.Worksheets.Item("Heat Generators")

# Update renewable factor (column K) values to 0 for specified rows
$ws1.Range("K4").Value = 0
$ws1.Range("K5").Value = 0
$ws1.Range("K6").Value = 0
$ws1.Range("K7").Value = 0
$ws1.Range("K8").Value = 0
$ws1.Range("K9").Value = 0
$ws1.Range("K13").Value = 0
$ws1.Range("K14").Value = 0
$ws1.Range("K15").Value = 0
$ws1.Range("K16").Value = 0
$ws1.Range("K17").Value = 0
$ws1.Range("K22").Value = 0
$ws1.Range("K23").Value = 0
$ws1.Range("K24").Value = 0
$ws1.Range("K25").Value = 0
$ws1.Range("K26").Value = 0
$ws1.Range("K27").Value = 0
$ws1.Range("K28").Value = 0
$ws1.Range("K29").Value = 0

$ws1.Activate()
$ws1.Range("I35").Select()

# --- Sheet: financal and other parameteres ---
$ws3 = $wb.Worksheets.Item("financal and other parameteres")
$ws3.Range("B3").Value = 0.05

$ws3.Activate()
$ws3.Range("D17").Select()

# --- Sheet: Heat Storage ---
$ws4 = $wb.Worksheets.Item("Heat Storage")
$ws4.Activate()
$ws4.Range("D15:D16").Select()

# Re-activate the originally active sheet (Heat Generators, tabSelected=1)
$ws1.Activate()
